$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. The existing "总计" sheet becomes "2022-Q1" (keeps its underlying sheet
#    identity / position-3 slot), and a brand-new "总计" sheet is appended
#    right after it to hold the refreshed summary table.
# ---------------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

$total = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q1)
$total.Name = "总计"

# Grab a reference cell that already carries the workbook's "header / index"
# style (s="2") so we can clone it onto the new sheets instead of reinventing
# the formatting.
$styleSrc = $wb.Worksheets.Item("2021-Q4").Range("B1")

# ---------------------------------------------------------------------------
# 2. Rebuild "2022-Q1" (the fund-holding detail table) from scratch.
# ---------------------------------------------------------------------------
$q1.Cells.Clear()

$q1Headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $q1Headers.Length; $i++) {
    $q1.Cells.Item(1, 2 + $i).Value = $q1Headers[$i]
}

$q1Rows = @(
    @("014185", "招商专精特新股票A",               "8.37", "30.94", "2.22", "0.1858", 3),
    @("000264", "博时内需增长混合",                 "3.87", "75.26", "4.64", "0.1796", 7),
    @("050012", "博时策略混合",                     "3.00", "73.86", "5.13", "0.1539", 5),
    @("014186", "招商专精特新股票C",                "3.46", "30.94", "2.22", "0.0768", 3),
    @("012153", "博时研究慧选混合型证券投资基金A",   "1.63", "75.28", "2.94", "0.0479", 10),
    @("005459", "银河嘉谊灵活配置混合A",             "6.47", "39.69", "0.62", "0.0401", 6),
    @("004677", "博时战略新兴产业混合",               "0.41", "89.27", "4.88", "0.0200", 6),
    @("005460", "银河嘉谊灵活配置混合C",             "2.79", "39.69", "0.62", "0.0173", 6),
    @("012154", "博时研究慧选混合型证券投资基金C",   "0.21", "75.28", "2.94", "0.0062", 10)
)

# Columns B, D, E, F, G hold numeric-looking values that must stay TEXT
# (e.g. "014185", "8.37", "0.1858") — force text entry, then snap the style
# back to Normal so only the cell's stored type (not its formatting) differs
# from a plain, unstyled cell. Column C (fund name) is plain text already
# and column H (ranking) is a genuine number, so neither needs this dance.
for ($r = 0; $r -lt $q1Rows.Length; $r++) {
    $row = 2 + $r
    $data = $q1Rows[$r]

    $q1.Cells.Item($row, 1).Value = $r

    $q1.Range($q1.Cells.Item($row, 2), $q1.Cells.Item($row, 2)).NumberFormat = "@"
    $q1.Cells.Item($row, 2).Value = $data[0]

    $q1.Cells.Item($row, 3).Value = $data[1]

    $q1.Range($q1.Cells.Item($row, 4), $q1.Cells.Item($row, 7)).NumberFormat = "@"
    $q1.Cells.Item($row, 4).Value = $data[2]
    $q1.Cells.Item($row, 5).Value = $data[3]
    $q1.Cells.Item($row, 6).Value = $data[4]
    $q1.Cells.Item($row, 7).Value = $data[5]

    $q1.Range($q1.Cells.Item($row, 2), $q1.Cells.Item($row, 2)).Style = "Normal"
    $q1.Range($q1.Cells.Item($row, 4), $q1.Cells.Item($row, 7)).Style = "Normal"

    $q1.Cells.Item($row, 8).Value = $data[6]
}

# Apply the shared "header / index" style to row 1 (B:H) and column A
# (rows 2-10) the same way the sibling quarterly sheets do it.
$styleSrc.Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$styleSrc.Copy()
$q1.Range("A2:A10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Populate the new "总计" sheet with the refreshed summary table
#    (2022-Q1 on top, followed by the previously existing quarters).
# ---------------------------------------------------------------------------
$totalHeaders = @("日期", "持有数量(只)", "持有市值(亿元)")
for ($i = 0; $i -lt $totalHeaders.Length; $i++) {
    $total.Cells.Item(1, 2 + $i).Value = $totalHeaders[$i]
}

$totalRows = @(
    @("2022-Q1", 9, 0.73),
    @("2021-Q4", 6, 0.35),
    @("2021-Q3", 3, 0.45)
)

for ($r = 0; $r -lt $totalRows.Length; $r++) {
    $row = 2 + $r
    $data = $totalRows[$r]

    $total.Cells.Item($row, 1).Value = $r

    $total.Range($total.Cells.Item($row, 2), $total.Cells.Item($row, 2)).NumberFormat = "@"
    $total.Cells.Item($row, 2).Value = $data[0]
    $total.Range($total.Cells.Item($row, 2), $total.Cells.Item($row, 2)).Style = "Normal"

    $total.Cells.Item($row, 3).Value = $data[1]
    $total.Cells.Item($row, 4).Value = $data[2]
}

$styleSrc.Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$styleSrc.Copy()
$total.Range("A2:A4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$q1.Range("A1").Select()
